$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2835.2173
$ws.Range("I86").Value = 2119.125
$ws.Range("J86").Value = 4472
$ws.Range("K86").Value = 2119.125
$ws.Range("L86").Value = 4472
$ws.Range("M86").Value = -996.125
$ws.Range("N86").Value = -6718
$ws.Range("H89").Value = 2835.2173
$ws.Range("I89").Value = 2119.125
$ws.Range("J89").Value = 4472
$ws.Range("K89").Value = 10595.625
$ws.Range("L89").Value = 22360
$ws.Range("M89").Value = -4979.625
$ws.Range("N89").Value = -33592
$ws.Range("H99").Value = 1059.8
$ws.Range("I99").Value = 800
$ws.Range("J99").Value = 1449.5
$ws.Range("K99").Value = 2400
$ws.Range("L99").Value = 4348.5
$ws.Range("M99").Value = -902
$ws.Range("N99").Value = -7344.5
$ws.Range("H106").Value = 1829.8889
$ws.Range("I106").Value = 1829.8889
$ws.Range("K106").Value = 1829.8889
$ws.Range("M106").Value = -1198.8889
$ws.Range("H112").Value = 41460
$ws.Range("J112").Value = 42996.92
$ws.Range("L112").Value = 128990.76
$ws.Range("N112").Value = -131206.76
$ws.Range("H135").Value = 1176.4482
$ws.Range("I135").Value = 1053.1666
$ws.Range("K135").Value = 9478.499400000001
$ws.Range("M135").Value = -6943.499400000001
$ws.Range("H137").Value = 1721.0975
$ws.Range("J137").Value = 1612.4286
$ws.Range("L137").Value = 4837.2858
$ws.Range("N137").Value = -9937.2858
$ws.Range("H138").Value = 8336394
$ws.Range("I138").Value = 1425.7693
$ws.Range("J138").Value = 10641811
$ws.Range("K138").Value = 4277.3079
$ws.Range("L138").Value = 31925433
$ws.Range("M138").Value = 862.6921000000002
$ws.Range("N138").Value = -31935713

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 63672.895
$ws.Range("I5").Value = 93000.38
$ws.Range("K5").Value = 93000.38
$ws.Range("M5").Value = -92888.38
$ws.Range("H44").Value = 34999.5
$ws.Range("J44").Value = 34999.5
$ws.Range("L44").Value = 34999.5
$ws.Range("N44").Value = -35975.5
$ws.Range("H45").Value = 7173.8887
$ws.Range("I45").Value = 8684.429
$ws.Range("K45").Value = 8684.429
$ws.Range("M45").Value = -8307.429
$ws.Range("H61").Value = 174683.11
$ws.Range("I61").Value = 2500.7727
$ws.Range("K61").Value = 2500.7727
$ws.Range("M61").Value = -2288.7727
$ws.Range("H74").Value = 5084.7925
$ws.Range("I74").Value = 1183.8889
$ws.Range("J74").Value = 27027.375
$ws.Range("K74").Value = 1183.8889
$ws.Range("L74").Value = 27027.375
$ws.Range("M74").Value = -309.8888999999999
$ws.Range("N74").Value = -28775.375
$ws.Range("H77").Value = 5084.7925
$ws.Range("I77").Value = 1183.8889
$ws.Range("J77").Value = 27027.375
$ws.Range("K77").Value = 5919.4445
$ws.Range("L77").Value = 135136.875
$ws.Range("M77").Value = -1551.4445
$ws.Range("N77").Value = -143872.875
$ws.Range("H97").Value = 1127.3846
$ws.Range("I97").Value = 855.5278
$ws.Range("K97").Value = 855.5278
$ws.Range("M97").Value = -359.5278
$ws.Range("H122").Value = 1476.8334
$ws.Range("I122").Value = 1172.2
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 3516.6
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -1066.6
$ws.Range("N122").Value = -13900
$ws.Range("H136").Value = 174683.11
$ws.Range("I136").Value = 2500.7727
$ws.Range("K136").Value = 7502.3181
$ws.Range("M136").Value = -4952.3181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 63672.895
$ws.Range("I4").Value = 93000.38
$ws.Range("K4").Value = 93000.38
$ws.Range("M4").Value = -92885.38
$ws.Range("H94").Value = 915.0789
$ws.Range("I94").Value = 618.89655
$ws.Range("J94").Value = 1869.4445
$ws.Range("K94").Value = 618.89655
$ws.Range("L94").Value = 1869.4445
$ws.Range("M94").Value = -167.89655
$ws.Range("N94").Value = -2771.4445

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 47501.637
$ws.Range("J31").Value = 2234.4167
$ws.Range("L31").Value = 2234.4167
$ws.Range("N31").Value = -2824.4167
$ws.Range("H34").Value = 47501.637
$ws.Range("J34").Value = 2234.4167
$ws.Range("L34").Value = 2234.4167
$ws.Range("N34").Value = -2638.4167
$ws.Range("H53").Value = 7633.3335
$ws.Range("J53").Value = 7633.3335
$ws.Range("L53").Value = 7633.3335
$ws.Range("N53").Value = -8847.333500000001
$ws.Range("H58").Value = 2303.8147
$ws.Range("J58").Value = 3388
$ws.Range("L58").Value = 3388
$ws.Range("N58").Value = -3794
$ws.Range("H107").Value = 985.44446
$ws.Range("I107").Value = 586.6667
$ws.Range("K107").Value = 586.6667
$ws.Range("M107").Value = 1333.3333
$ws.Range("H108").Value = 35750
$ws.Range("J108").Value = 35750
$ws.Range("L108").Value = 35750
$ws.Range("N108").Value = -43430
$ws.Range("H132").Value = 3323.44
$ws.Range("I132").Value = 2742.913
$ws.Range("K132").Value = 8228.739
$ws.Range("M132").Value = -5698.739
$ws.Range("H136").Value = 2303.8147
$ws.Range("J136").Value = 3388
$ws.Range("L136").Value = 10164
$ws.Range("N136").Value = -15264
$ws.Range("H141").Value = 142421
$ws.Range("J141").Value = 142421
$ws.Range("L141").Value = 142421
$ws.Range("N141").Value = -152781

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 786.5294
$ws.Range("I5").Value = 786.5294
$ws.Range("K5").Value = 2359.5882
$ws.Range("M5").Value = -2247.5882
$ws.Range("H26").Value = 1719.9688
$ws.Range("I26").Value = 1068.091
$ws.Range("J26").Value = 2061.4285
$ws.Range("K26").Value = 3204.273
$ws.Range("L26").Value = 6184.2855
$ws.Range("M26").Value = -2916.273
$ws.Range("N26").Value = -6760.2855
$ws.Range("H131").Value = 25639.273
$ws.Range("J131").Value = 3686.7812
$ws.Range("L131").Value = 11060.3436
$ws.Range("N131").Value = -21140.3436
$ws.Range("H135").Value = 786.5294
$ws.Range("I135").Value = 786.5294
$ws.Range("K135").Value = 7078.7646
$ws.Range("M135").Value = -4543.7646

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 19124.25
$ws.Range("J70").Value = 32249.75
$ws.Range("L70").Value = 32249.75
$ws.Range("N70").Value = -32789.75
$ws.Range("H73").Value = 19124.25
$ws.Range("J73").Value = 32249.75
$ws.Range("L73").Value = 32249.75
$ws.Range("N73").Value = -34121.75
$ws.Range("H122").Value = 2289.35
$ws.Range("J122").Value = 2353.5
$ws.Range("L122").Value = 7060.5
$ws.Range("N122").Value = -11960.5
$ws.Range("H126").Value = 14004.95
$ws.Range("I126").Value = 17393.6
$ws.Range("J126").Value = 3839
$ws.Range("K126").Value = 52180.8
$ws.Range("L126").Value = 11517
$ws.Range("M126").Value = -49710.8
$ws.Range("N126").Value = -16457
$ws.Range("H136").Value = 32237.96
$ws.Range("J136").Value = 32237.96
$ws.Range("L136").Value = 96713.88
$ws.Range("N136").Value = -101813.88

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3461.543
$ws.Range("I7").Value = 2229.3333
$ws.Range("J7").Value = 6150
$ws.Range("K7").Value = 2229.3333
$ws.Range("L7").Value = 6150
$ws.Range("M7").Value = -2117.3333
$ws.Range("N7").Value = -6374
$ws.Range("H40").Value = 3523.0952
$ws.Range("I40").Value = 2695.3572
$ws.Range("K40").Value = 2695.3572
$ws.Range("M40").Value = -2559.3572
$ws.Range("H122").Value = 4535.6787
$ws.Range("I122").Value = 4531.3076
$ws.Range("K122").Value = 13593.9228
$ws.Range("M122").Value = -11143.9228
$ws.Range("H126").Value = 3461.543
$ws.Range("I126").Value = 2229.3333
$ws.Range("J126").Value = 6150
$ws.Range("K126").Value = 6687.999899999999
$ws.Range("L126").Value = 18450
$ws.Range("M126").Value = -4217.999899999999
$ws.Range("N126").Value = -23390

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3713.3333
$ws.Range("I62").Value = 2430.8333
$ws.Range("J62").Value = 4995.8335
$ws.Range("K62").Value = 2430.8333
$ws.Range("L62").Value = 4995.8335
$ws.Range("M62").Value = -1806.8333
$ws.Range("N62").Value = -6243.8335
$ws.Range("H65").Value = 3713.3333
$ws.Range("I65").Value = 2430.8333
$ws.Range("J65").Value = 4995.8335
$ws.Range("K65").Value = 12154.1665
$ws.Range("L65").Value = 24979.1675
$ws.Range("M65").Value = -9034.166499999999
$ws.Range("N65").Value = -31219.1675
$ws.Range("H81").Value = 7532.5
$ws.Range("I81").Value = 10338
$ws.Range("K81").Value = 20676
$ws.Range("M81").Value = -19615
$ws.Range("H84").Value = 7532.5
$ws.Range("I84").Value = 10338
$ws.Range("K84").Value = 103380
$ws.Range("M84").Value = -98076
